# Refresh the cryptos price table (columns D "Price" and E "Volume(1h)")
# plus a Mantle/Filecoin row re-order, per the automated GitHub Actions
# "Updated cryptos list" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.720.82'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '3.783.84'
$ws.Range('E3').Value = '  -1.70%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '''597.61'
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').Value = '''169.01'
$ws.Range('E6').Value = '  +0.79%  '
$ws.Range('D7').Value = '3.782.81'
$ws.Range('E7').Value = '  -1.56%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').Value = '''0.530'
$ws.Range('E9').Value = '  +0.57%  '
$ws.Range('D10').Value = '''0.167'
$ws.Range('E10').Value = '  +2.34%  '
$ws.Range('D11').Value = '''6.51'
$ws.Range('E11').Value = '  +1.97%  '
$ws.Range('D12').Value = '''0.460'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').Value = '''0.0000273'
$ws.Range('E13').Value = '  +6.65%  '
$ws.Range('D14').Value = '''36.92'
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').Value = '4.420.13'
$ws.Range('E15').Value = '  -1.71%  '
$ws.Range('D16').Value = '3.788.20'
$ws.Range('E16').Value = '  -1.76%  '
$ws.Range('D17').Value = '''19.03'
$ws.Range('E17').Value = '  +5.64%  '
$ws.Range('D18').Value = '67.764.02'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').Value = '''7.27'
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('D21').Value = '''10.59'
$ws.Range('E21').Value = '  -1.33%  '
$ws.Range('D22').Value = '''467.03'
$ws.Range('E22').Value = '  +0.68%  '
$ws.Range('D23').Value = '''0.727'
$ws.Range('E23').Value = '  -0.66%  '
$ws.Range('E24').Value = '  -5.20%  '
$ws.Range('D25').Value = '''83.45'
$ws.Range('E25').Value = '  +0.74%  '
$ws.Range('E26').Value = '  +1.86%  '
$ws.Range('D27').Value = '''12.20'
$ws.Range('E27').Value = '  +1.78%  '
$ws.Range('D28').Value = '''10.32'
$ws.Range('E28').Value = '  +4.16%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = '''2.93'
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('D31').Value = '3.939.34'
$ws.Range('E31').Value = '  -1.58%  '
$ws.Range('D32').Value = '''7.63'
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').Value = '''2.26'
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('D34').Value = '''30.49'
$ws.Range('E34').Value = '  -1.56%  '
$ws.Range('D35').Value = '''9.17'
$ws.Range('E35').Value = '  -2.52%  '
$ws.Range('D36').Value = '3.751.36'
$ws.Range('E36').Value = '  -1.65%  '
$ws.Range('D37').Value = '''3.81'
$ws.Range('E37').Value = '  +4.20%  '
$ws.Range('E38').Value = '  +1.45%  '
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').Value = '''5.91'
$ws.Range('E39').Value = '  +0.77%  '
$ws.Range('B40').Value = 'Mantle'
$ws.Range('C40').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D40').Value = '''1.01'
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('E41').Value = '  -1.01%  '
$ws.Range('D42').Value = '''0.999'
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '''0.318'
$ws.Range('E43').Value = '  +2.65%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E45').Value = '  +2.05%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('D47').Value = '''406.48'
$ws.Range('E47').Value = '  -3.81%  '
$ws.Range('D48').Value = '''46.28'
$ws.Range('E48').Value = '  -1.63%  '
$ws.Range('D49').Value = '''0.000278'
$ws.Range('E49').Value = '  -5.47%  '
$ws.Range('D50').Value = '''142.10'
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('E51').Value = '  +0.25%  '
